$d = $word.ActiveDocument

# --- Change 1: update the DATE field text, preserving <w:noProof/> formatting ---
$dateRange = $d.Content
$null = $dateRange.Find.Execute("Wednesday, 15 April 2009")
$dateRange.NoProofing = $true
$dateRange.Text = "Thursday, 16 April 2009"

# --- Change 2: insert the new "User interface" / "Project explorer" sections ---
$count = $d.Paragraphs.Count
$anchor = $d.Paragraphs.Item($count - 1)
$anchorRange = $anchor.Range
$null = $anchorRange.InsertParagraphAfter()
$target = $d.Paragraphs.Item($count)
$xmlBlock = '<w:p><w:r><w:t>User interface:</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t>Provide easy access to the most current operations. Provide toolbars for the setting of the physical data.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Allow marking of geometry elements and regions with capabilities. These can then be translated by the system into numerical data. E.g. an element is </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>rigid,</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> a curve is a wall or a symmetry line.</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>Project explorer:</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t>The project explorer is used to browse through existing projects / experiments and create new ones. The project explorer cannot be used to visualize any of the data in the project / experiment but it can be used to create new child-experiments etc.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t>The project explorer should allow storing a script of the changes to be made so that the changes can be repeated easily over and over again.</w:t></w:r></w:p>'
$null = $target.Range.InsertXML($xmlBlock)
